# The needle calibration data (rows 2-18, columns A-D) needs to be
# re-sorted in ascending order of column A ("time (s)"). The header
# row (row 1) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:D18")
$keyRange = $ws.Range("A2:A18")

$dataRange.Sort($keyRange, 1)
